# Apply cryptos-list refresh: updated prices / % changes, and the
# EthereumClassic <-> Monero row swap (rows 31-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.582.07"
$ws.Range('E2').Value = "  -1.43%  "
$ws.Range('D3').Value = "'2.288.99"
$ws.Range('E3').Value = "  +1.10%  "
$ws.Range('E4').Value = "  -0.04%  "
$ws.Range('D5').Value = "'95.56"
$ws.Range('E5').Value = "  -3.45%  "
$ws.Range('D6').Value = "'267.93"
$ws.Range('E6').Value = "  -2.97%  "
$ws.Range('E7').Value = "  -0.94%  "
$ws.Range('D8').Value = "'0.999"
$ws.Range('D9').Value = "'0.604"
$ws.Range('E9').Value = "  -5.22%  "
$ws.Range('D10').Value = "'44.75"
$ws.Range('E10').Value = "  -7.62%  "
$ws.Range('D11').Value = "'0.0936"
$ws.Range('E11').Value = "  -1.00%  "
$ws.Range('D12').Value = "'7.81"
$ws.Range('E12').Value = "  -5.12%  "
$ws.Range('D13').Value = "'0.105"
$ws.Range('E13').Value = "  +0.35%  "
$ws.Range('D14').Value = "'2.631.83"
$ws.Range('E14').Value = "  +1.28%  "
$ws.Range('E15').Value = "  -2.69%  "
$ws.Range('D16').Value = "'0.844"
$ws.Range('E16').Value = "  +0.13%  "
$ws.Range('D17').Value = "'2.289.40"
$ws.Range('E17').Value = "  +1.52%  "
$ws.Range('D18').Value = "'43.553.95"
$ws.Range('E18').Value = "  -1.52%  "
$ws.Range('E19').Value = "  +0.81%  "
$ws.Range('E20').Value = "  -0.92%  "
$ws.Range('D21').Value = "'72.43"
$ws.Range('E21').Value = "  +1.90%  "
$ws.Range('E22').Value = "  +5.59%  "
$ws.Range('D23').Value = "'235.17"
$ws.Range('E23').Value = "  -0.15%  "
$ws.Range('D24').Value = "'9.08"
$ws.Range('E24').Value = "  -15.48%  "
$ws.Range('E26').Value = "  -1.95%  "
$ws.Range('D27').Value = "'11.20"
$ws.Range('E27').Value = "  -3.17%  "
$ws.Range('E28').Value = "  +2.53%  "
$ws.Range('D29').Value = "'40.07"
$ws.Range('E29').Value = "  -0.02%  "
$ws.Range('E30').Value = "  -2.02%  "
$ws.Range('B31').Value = "Monero"
$ws.Range('C31').Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D31').Value = "'174.80"
$ws.Range('E31').Value = "  +0.83%  "
$ws.Range('B32').Value = "EthereumClassic"
$ws.Range('C32').Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('D32').Value = "'22.10"
$ws.Range('E32').Value = "  +3.73%  "
$ws.Range('D33').Value = "'0.0882"
$ws.Range('E33').Value = "  -4.65%  "
$ws.Range('E34').Value = "  -7.37%  "
$ws.Range('E35').Value = "  -0.10%  "
$ws.Range('E36').Value = "  -5.13%  "
$ws.Range('D37').Value = "'0.0357"
$ws.Range('E37').Value = "  +0.49%  "
$ws.Range('D38').Value = "'4.38"
$ws.Range('E38').Value = "  -1.35%  "
$ws.Range('D39').Value = "'3.30"
$ws.Range('E39').Value = "  -7.64%  "
$ws.Range('E40').Value = "  -7.59%  "
$ws.Range('E41').Value = "  +5.79%  "
$ws.Range('D42').Value = "'65.17"
$ws.Range('E42').Value = "  +3.36%  "
$ws.Range('D43').Value = "'12.10"
$ws.Range('E43').Value = "  -5.10%  "
$ws.Range('D45').Value = "'8.80"
$ws.Range('E45').Value = "  +2.50%  "
$ws.Range('E46').Value = "  -4.91%  "
$ws.Range('E47').Value = "  -1.32%  "
$ws.Range('D48').Value = "'98.19"
$ws.Range('E48').Value = "  -2.48%  "
$ws.Range('E49').Value = "  -0.64%  "
$ws.Range('E50').Value = "  +6.57%  "
$ws.Range('D51').Value = "'2.512.08"
$ws.Range('E51').Value = "  +1.23%  "
